$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right 105 -> 84, Wrong -2 -> -4, Max text "105 / 140" -> "80 / 112"
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "80 / 112"
